$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The author's commit removed the single line reading "AHMED ESSAM"
# (an entire run of text) from its own paragraph, and - because the
# paragraph no longer has any visible text - Word also dropped the
# now-superfluous centered-justification (<w:jc w:val="center"/>)
# paragraph property, leaving an empty, left-aligned paragraph behind.
# ------------------------------------------------------------------

$needle = "AHMED ESSAM"

# Step 1: locate the paragraph that currently holds the target text so
# we can fix up its paragraph-level formatting afterwards, even though
# the text itself is about to disappear (and with it, our ability to
# Find() it again).
$locateRange = $d.Content
$locateFind = $locateRange.Find
$locateFind.ClearFormatting()
$found = $locateFind.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$targetParaIndex = -1
if ($found) {
    $paragraphs = $d.Paragraphs
    for ($i = 1; $i -le $paragraphs.Count; $i++) {
        $candidate = $paragraphs.Item($i)
        if ($candidate.Range.Start -le $locateRange.Start -and $candidate.Range.End -ge $locateRange.End) {
            $targetParaIndex = $i
            break
        }
    }
}

# Step 2: delete the text. Using Find/Replace (rather than just blanking
# the Range.Text) also removes the now-empty <w:r> run entirely, which
# matches how the document was actually edited.
$deleteFind = $d.Content.Find
$deleteFind.ClearFormatting()
$deleteFind.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# Step 3: the paragraph is now empty; remove its centered alignment so it
# reverts to the default (left) alignment, matching the target markup
# (no <w:jc> element at all).
if ($targetParaIndex -gt 0) {
    $targetPara = $d.Paragraphs.Item($targetParaIndex)
    if ($targetPara.Range.Text -eq "" -or $targetPara.Range.Text -eq [char]13) {
        $targetPara.Alignment = 0
    }
}
